$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 108.333336
$ws.Cells.Item(12, 9).Value = 108.333336
$ws.Cells.Item(12, 11).Value = 108.333336
$ws.Cells.Item(12, 13).Value = 61.666664
$ws.Cells.Item(32, 8).Value = 2859.818
$ws.Cells.Item(32, 9).Value = 1699.5
$ws.Cells.Item(32, 10).Value = 3117.6667
$ws.Cells.Item(32, 11).Value = 1699.5
$ws.Cells.Item(32, 12).Value = 3117.6667
$ws.Cells.Item(32, 13).Value = -1373.5
$ws.Cells.Item(32, 14).Value = -3769.6667
$ws.Cells.Item(33, 8).Value = 74.72727
$ws.Cells.Item(33, 9).Value = 76.388885
$ws.Cells.Item(33, 11).Value = 76.388885
$ws.Cells.Item(33, 13).Value = 152.611115
$ws.Cells.Item(107, 8).Value = 1646.8334
$ws.Cells.Item(107, 9).Value = 424.66666
$ws.Cells.Item(107, 11).Value = 424.66666
$ws.Cells.Item(107, 13).Value = 1495.33334
$ws.Cells.Item(113, 8).Value = 6663.1113
$ws.Cells.Item(113, 9).Value = 5967.4443
$ws.Cells.Item(113, 10).Value = 8054.4443
$ws.Cells.Item(113, 11).Value = 5967.4443
$ws.Cells.Item(113, 12).Value = 8054.4443
$ws.Cells.Item(113, 13).Value = -2713.4443
$ws.Cells.Item(113, 14).Value = -14562.4443
$ws.Cells.Item(116, 8).Value = 23749.4
$ws.Cells.Item(116, 9).Value = 8719.75
$ws.Cells.Item(116, 11).Value = 8719.75
$ws.Cells.Item(116, 13).Value = -5277.75
$ws.Cells.Item(137, 8).Value = 43487580
$ws.Cells.Item(137, 9).Value = 142860720
$ws.Cells.Item(137, 10).Value = 11835.625
$ws.Cells.Item(137, 11).Value = 428582160
$ws.Cells.Item(137, 12).Value = 35506.875
$ws.Cells.Item(137, 13).Value = -428579610
$ws.Cells.Item(137, 14).Value = -40606.875
$ws.Cells.Item(138, 8).Value = 4934.2183
$ws.Cells.Item(138, 10).Value = 5155.512
$ws.Cells.Item(138, 12).Value = 15466.536
$ws.Cells.Item(138, 14).Value = -25746.536

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2202.9
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 14).Value = ""
$ws.Cells.Item(32, 8).Value = 172448.64
$ws.Cells.Item(32, 9).Value = 198592.64
$ws.Cells.Item(32, 10).Value = 36499.9
$ws.Cells.Item(32, 11).Value = 198592.64
$ws.Cells.Item(32, 12).Value = 36499.9
$ws.Cells.Item(32, 13).Value = -198305.64
$ws.Cells.Item(32, 14).Value = -37073.9
$ws.Cells.Item(45, 8).Value = 1512
$ws.Cells.Item(45, 9).Value = 1186.6666
$ws.Cells.Item(45, 11).Value = 1186.6666
$ws.Cells.Item(45, 13).Value = -809.6666
$ws.Cells.Item(110, 8).Value = 1083.8667
$ws.Cells.Item(110, 9).Value = 734.7727
$ws.Cells.Item(110, 10).Value = 2043.875
$ws.Cells.Item(110, 11).Value = 734.7727
$ws.Cells.Item(110, 12).Value = 2043.875
$ws.Cells.Item(110, 13).Value = 1310.2273
$ws.Cells.Item(110, 14).Value = -6133.875
$ws.Cells.Item(116, 8).Value = 2202.9
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).Value = ""
$ws.Cells.Item(132, 8).Value = 1044041.25
$ws.Cells.Item(132, 9).Value = 1564376.4
$ws.Cells.Item(132, 11).Value = 4693129.199999999
$ws.Cells.Item(132, 13).Value = -4690599.199999999
$ws.Cells.Item(134, 8).Value = 74932.336
$ws.Cells.Item(134, 10).Value = 74932.336
$ws.Cells.Item(134, 12).Value = 74932.336
$ws.Cells.Item(134, 14).Value = -85072.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2202.9
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 14).Value = ""
$ws.Cells.Item(35, 8).Value = 34500
$ws.Cells.Item(35, 10).Value = 34500
$ws.Cells.Item(35, 12).Value = 34500
$ws.Cells.Item(35, 14).Value = -35120
$ws.Cells.Item(107, 8).Value = 16667614
$ws.Cells.Item(107, 9).Value = 31250688
$ws.Cells.Item(107, 10).Value = 1242.7142
$ws.Cells.Item(107, 11).Value = 31250688
$ws.Cells.Item(107, 12).Value = 1242.7142
$ws.Cells.Item(107, 13).Value = -31248768
$ws.Cells.Item(107, 14).Value = -5082.7142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 85131.25
$ws.Cells.Item(22, 9).Value = 270.25
$ws.Cells.Item(22, 10).Value = 127561.75
$ws.Cells.Item(22, 11).Value = 270.25
$ws.Cells.Item(22, 12).Value = 127561.75
$ws.Cells.Item(22, 13).Value = 79.75
$ws.Cells.Item(22, 14).Value = -128261.75
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 14).Value = ""
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 14).Value = ""
$ws.Cells.Item(105, 8).Value = 5184.3076
$ws.Cells.Item(105, 9).Value = 5371.4287
$ws.Cells.Item(105, 11).Value = 5371.4287
$ws.Cells.Item(105, 13).Value = -3624.4287
$ws.Cells.Item(107, 8).Value = 382.72726
$ws.Cells.Item(107, 9).Value = 351.2353
$ws.Cells.Item(107, 11).Value = 351.2353
$ws.Cells.Item(107, 13).Value = 1568.7647
$ws.Cells.Item(122, 8).Value = 19634.834
$ws.Cells.Item(122, 9).Value = 7803
$ws.Cells.Item(122, 10).Value = 36199.4
$ws.Cells.Item(122, 11).Value = 23409
$ws.Cells.Item(122, 12).Value = 108598.2
$ws.Cells.Item(122, 13).Value = -20959
$ws.Cells.Item(122, 14).Value = -113498.2
$ws.Cells.Item(134, 8).Value = 3404.762
$ws.Cells.Item(134, 9).Value = 2901.5
$ws.Cells.Item(134, 11).Value = 8704.5
$ws.Cells.Item(134, 13).Value = -6169.5
$ws.Cells.Item(141, 8).Value = 203743
$ws.Cells.Item(141, 10).Value = 215784
$ws.Cells.Item(141, 12).Value = 215784
$ws.Cells.Item(141, 14).Value = -226144

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 178.55556
$ws.Cells.Item(2, 9).Value = 71.7
$ws.Cells.Item(2, 10).Value = 312.125
$ws.Cells.Item(2, 11).Value = 430.2
$ws.Cells.Item(2, 12).Value = 1872.75
$ws.Cells.Item(2, 13).Value = -317.2
$ws.Cells.Item(2, 14).Value = -2098.75
$ws.Cells.Item(17, 8).Value = 2309.25
$ws.Cells.Item(17, 9).Value = 1594.8
$ws.Cells.Item(17, 10).Value = 3500
$ws.Cells.Item(17, 11).Value = 4784.4
$ws.Cells.Item(17, 12).Value = 10500
$ws.Cells.Item(17, 13).Value = -4615.4
$ws.Cells.Item(17, 14).Value = -10838
$ws.Cells.Item(34, 8).Value = 10.8
$ws.Cells.Item(34, 9).Value = 10.8
$ws.Cells.Item(34, 11).Value = 32.40000000000001
$ws.Cells.Item(34, 13).Value = 51.59999999999999
$ws.Cells.Item(38, 8).Value = 74.125
$ws.Cells.Item(38, 10).Value = 153.71428
$ws.Cells.Item(38, 12).Value = 461.14284
$ws.Cells.Item(38, 14).Value = -1155.14284
$ws.Cells.Item(39, 8).Value = 2066.6667
$ws.Cells.Item(39, 10).Value = 5000
$ws.Cells.Item(39, 12).Value = 15000
$ws.Cells.Item(39, 14).Value = -15588
$ws.Cells.Item(70, 8).Value = 1953.6
$ws.Cells.Item(70, 9).Value = 427
$ws.Cells.Item(70, 10).Value = 2971.3333
$ws.Cells.Item(70, 11).Value = 1281
$ws.Cells.Item(70, 12).Value = 8913.999899999999
$ws.Cells.Item(70, 13).Value = -966
$ws.Cells.Item(70, 14).Value = -9543.999899999999
$ws.Cells.Item(73, 8).Value = 1953.6
$ws.Cells.Item(73, 9).Value = 427
$ws.Cells.Item(73, 10).Value = 2971.3333
$ws.Cells.Item(73, 11).Value = 1281
$ws.Cells.Item(73, 12).Value = 8913.999899999999
$ws.Cells.Item(73, 13).Value = -189
$ws.Cells.Item(73, 14).Value = -11097.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 8117.3887
$ws.Cells.Item(57, 10).Value = 22900
$ws.Cells.Item(57, 12).Value = 22900
$ws.Cells.Item(57, 14).Value = -24540
$ws.Cells.Item(102, 8).Value = 926.2941
$ws.Cells.Item(102, 9).Value = 723.625
$ws.Cells.Item(102, 11).Value = 723.625
$ws.Cells.Item(102, 13).Value = 898.375
$ws.Cells.Item(122, 8).Value = 69976.125
$ws.Cells.Item(122, 9).Value = 117459.89
$ws.Cells.Item(122, 10).Value = 8925.571
$ws.Cells.Item(122, 11).Value = 352379.67
$ws.Cells.Item(122, 12).Value = 26776.713
$ws.Cells.Item(122, 13).Value = -349929.67
$ws.Cells.Item(122, 14).Value = -31676.713
$ws.Cells.Item(132, 8).Value = 16328.765
$ws.Cells.Item(132, 9).Value = 12182.333
$ws.Cells.Item(132, 11).Value = 36546.999
$ws.Cells.Item(132, 13).Value = -34016.999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 673.3333
$ws.Cells.Item(9, 9).Value = 760
$ws.Cells.Item(9, 11).Value = 760
$ws.Cells.Item(9, 13).Value = -536
$ws.Cells.Item(16, 8).Value = 930.13336
$ws.Cells.Item(16, 9).Value = 534.7692
$ws.Cells.Item(16, 10).Value = 3500
$ws.Cells.Item(16, 11).Value = 534.7692
$ws.Cells.Item(16, 12).Value = 3500
$ws.Cells.Item(16, 13).Value = -364.7692
$ws.Cells.Item(16, 14).Value = -3840
$ws.Cells.Item(40, 8).Value = 3207.5557
$ws.Cells.Item(40, 9).Value = 3455.7334
$ws.Cells.Item(40, 10).Value = 1966.6666
$ws.Cells.Item(40, 11).Value = 3455.7334
$ws.Cells.Item(40, 12).Value = 1966.6666
$ws.Cells.Item(40, 13).Value = -3319.7334
$ws.Cells.Item(40, 14).Value = -2238.6666
$ws.Cells.Item(128, 8).Value = 54000
$ws.Cells.Item(128, 10).Value = 54000
$ws.Cells.Item(128, 12).Value = 54000
$ws.Cells.Item(128, 14).Value = -63960
$ws.Cells.Item(132, 8).Value = 8992466
$ws.Cells.Item(132, 9).Value = 16695365
$ws.Cells.Item(132, 11).Value = 50086095
$ws.Cells.Item(132, 13).Value = -50083565
$ws.Cells.Item(135, 8).Value = 46633
$ws.Cells.Item(135, 10).Value = 46633
$ws.Cells.Item(135, 12).Value = 46633
$ws.Cells.Item(135, 14).Value = -56773
$ws.Cells.Item(138, 8).Value = 100429
$ws.Cells.Item(138, 10).Value = 100429
$ws.Cells.Item(138, 12).Value = 100429
$ws.Cells.Item(138, 14).Value = -110709

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1512497.5
$ws.Cells.Item(4, 9).Value = 1512497.5
$ws.Cells.Item(4, 11).Value = 1512497.5
$ws.Cells.Item(4, 13).Value = -1512384.5
$ws.Cells.Item(96, 8).Value = 22223752
$ws.Cells.Item(96, 9).Value = 27779208
$ws.Cells.Item(96, 10).Value = 1933.3334
$ws.Cells.Item(96, 11).Value = 27779208
$ws.Cells.Item(96, 12).Value = 1933.3334
$ws.Cells.Item(96, 13).Value = -27777835
$ws.Cells.Item(96, 14).Value = -4679.3334
$ws.Cells.Item(107, 8).Value = 19348.053
$ws.Cells.Item(107, 9).Value = 2389
$ws.Cells.Item(107, 10).Value = 38191.445
$ws.Cells.Item(107, 11).Value = 7167
$ws.Cells.Item(107, 12).Value = 114574.335
$ws.Cells.Item(107, 13).Value = -5247
$ws.Cells.Item(107, 14).Value = -118414.335
$ws.Cells.Item(140, 8).Value = 126897
$ws.Cells.Item(140, 10).Value = 126897
$ws.Cells.Item(140, 12).Value = 126897
$ws.Cells.Item(140, 14).Value = -137257
